$d = $word.ActiveDocument

# 1. Remove the "_GoBack" bookmark from the first paragraph (title)
foreach ($bm in @($d.Bookmarks)) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}

# 2. Remove the es-PE language formatting on the title paragraph's run(s)
$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range
$titleRange.LanguageID = 0

# 3. Append a new sentence to the end of the second paragraph, after
#    "Hence, we can make big money savings for the benefit of the company."
$d.Content.Find.Execute("Hence, we can make big money savings for the benefit of the company.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Hence, we can make big money savings for the benefit of the company. I think my natural gas project will be outstanding.", 2)
